# "error solve ifrs list" - replace the bogus bulk-scale financial figures
# that were scraped into the IFRS company_list sheet with the corrected
# (much smaller) figures, and blank out the forecast years (2019/12(E),
# 2020/12(E), 2021/12(E) -> rows 7-9) whose estimated data was wrong.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 2014/12 (IFRS연결)
$ws.Range("D2").Value = 6306
$ws.Range("E2").Value = 208
$ws.Range("F2").Value = 208
$ws.Range("G2").Value = 253
$ws.Range("H2").Value = 65
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 48
$ws.Range("K2").Value = 10811
$ws.Range("L2").Value = 5466
$ws.Range("M2").Value = 5345
$ws.Range("N2").Value = 4857
$ws.Range("O2").Value = 488
$ws.Range("P2").Value = 224
$ws.Range("Q2").Value = -79
$ws.Range("R2").Value = -1165
$ws.Range("S2").Value = 611
$ws.Range("T2").Value = 871
$ws.Range("U2").Value = -950
$ws.Range("V2").Value = 4203
$ws.Range("W2").Value = 3.3
$ws.Range("X2").Value = 1.03
$ws.Range("Y2").Value = 0.36
$ws.Range("Z2").Value = 0.63
$ws.Range("AA2").Value = 102.27
$ws.Range("AB2").Value = 2688.95
$ws.Range("AC2").Value = 375
$ws.Range("AD2").Value = 323.7
$ws.Range("AE2").Value = 103592
$ws.Range("AF2").Value = 1.17
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 0.82
$ws.Range("AI2").Value = 273.24
$ws.Range("AJ2").Value = 4481591

# Row 3: 2015/12 (IFRS연결)
$ws.Range("D3").Value = 7047
$ws.Range("E3").Value = 644
$ws.Range("F3").Value = 644
$ws.Range("G3").Value = 941
$ws.Range("H3").Value = 666
$ws.Range("I3").Value = 479
$ws.Range("J3").Value = 187
$ws.Range("K3").Value = 12235
$ws.Range("L3").Value = 6209
$ws.Range("M3").Value = 6026
$ws.Range("N3").Value = 5216
$ws.Range("O3").Value = 810
$ws.Range("P3").Value = 233
$ws.Range("Q3").Value = 809
$ws.Range("R3").Value = -1069
$ws.Range("S3").Value = 111
$ws.Range("T3").Value = 969
$ws.Range("U3").Value = -159
$ws.Range("V3").Value = 4298
$ws.Range("W3").Value = 9.130000000000001
$ws.Range("X3").Value = 9.449999999999999
$ws.Range("Y3").Value = 9.51
$ws.Range("Z3").Value = 5.78
$ws.Range("AA3").Value = 103.03
$ws.Range("AB3").Value = 2764.41
$ws.Range("AC3").Value = 10011
$ws.Range("AD3").Value = 16.93
$ws.Range("AE3").Value = 107237
$ws.Range("AF3").Value = 1.58
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 0.59
$ws.Range("AI3").Value = 10.16
$ws.Range("AJ3").Value = 4666791

# Row 4: 2016/12 (IFRS연결) - column O (자본총계(비지배)) no longer present
$ws.Range("D4").Value = 7262
$ws.Range("E4").Value = 760
$ws.Range("F4").Value = 760
$ws.Range("G4").Value = 1978
$ws.Range("H4").Value = 1757
$ws.Range("I4").Value = 1491
$ws.Range("J4").Value = 266
$ws.Range("K4").Value = 14022
$ws.Range("L4").Value = 5264
$ws.Range("M4").Value = 8759
$ws.Range("N4").Value = 8759
$ws.Range("P4").Value = 303
$ws.Range("Q4").Value = 305
$ws.Range("R4").Value = -669
$ws.Range("S4").Value = 125
$ws.Range("T4").Value = 645
$ws.Range("U4").Value = -341
$ws.Range("V4").Value = 3613
$ws.Range("W4").Value = 10.47
$ws.Range("X4").Value = 24.19
$ws.Range("Y4").Value = 21.33
$ws.Range("Z4").Value = 13.38
$ws.Range("AA4").Value = 60.1
$ws.Range("AB4").Value = 3260.01
$ws.Range("AC4").Value = 28395
$ws.Range("AD4").Value = 5.51
$ws.Range("AE4").Value = 144762
$ws.Range("AF4").Value = 1.08
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 0.64
$ws.Range("AI4").Value = 4.06
$ws.Range("AJ4").Value = 6050313
$ws.Range("O4").ClearContents()

# Row 5: 2017/12 (IFRS연결) - column J (당기순이익(비지배)) no longer present
$ws.Range("D5").Value = 6903
$ws.Range("E5").Value = 534
$ws.Range("F5").Value = 534
$ws.Range("G5").Value = 471
$ws.Range("H5").Value = 325
$ws.Range("I5").Value = 325
$ws.Range("K5").Value = 14594
$ws.Range("L5").Value = 5452
$ws.Range("M5").Value = 9142
$ws.Range("N5").Value = 9128
$ws.Range("O5").Value = 15
$ws.Range("P5").Value = 307
$ws.Range("Q5").Value = 845
$ws.Range("R5").Value = -263
$ws.Range("S5").Value = 273
$ws.Range("T5").Value = 250
$ws.Range("U5").Value = 595
$ws.Range("V5").Value = 3985
$ws.Range("W5").Value = 7.73
$ws.Range("X5").Value = 4.71
$ws.Range("Y5").Value = 3.64
$ws.Range("Z5").Value = 2.27
$ws.Range("AA5").Value = 59.63
$ws.Range("AB5").Value = 3331.15
$ws.Range("AC5").Value = 5310
$ws.Range("AD5").Value = 24.11
$ws.Range("AE5").Value = 148700
$ws.Range("AF5").Value = 0.86
$ws.Range("AG5").Value = 1000
$ws.Range("AH5").Value = 0.78
$ws.Range("AI5").Value = 18.86
$ws.Range("AJ5").Value = 6138228
$ws.Range("J5").ClearContents()

# Row 6: 2018/12 (IFRS연결)
$ws.Range("D6").Value = 7075
$ws.Range("E6").Value = 505
$ws.Range("F6").Value = 505
$ws.Range("G6").Value = -1731
$ws.Range("H6").Value = -1639
$ws.Range("I6").Value = -1640
$ws.Range("K6").Value = 12419
$ws.Range("L6").Value = 5062
$ws.Range("M6").Value = 7358
$ws.Range("N6").Value = 7342
$ws.Range("P6").Value = 307
$ws.Range("Q6").Value = 583
$ws.Range("R6").Value = -751
$ws.Range("S6").Value = -561
$ws.Range("T6").Value = 251
$ws.Range("U6").Value = 332
$ws.Range("V6").Value = 3709
$ws.Range("W6").Value = 7.14
$ws.Range("X6").Value = -23.16
$ws.Range("Y6").Value = -19.91
$ws.Range("Z6").Value = -12.13
$ws.Range("AA6").Value = 68.8
$ws.Range("AB6").Value = 2774.75
$ws.Range("AC6").Value = -26714
$ws.Range("AD6").Value = -4.08
$ws.Range("AE6").Value = 119609
$ws.Range("AF6").Value = 0.91
$ws.Range("AG6").Value = 1000
$ws.Range("AH6").Value = 0.92
$ws.Range("AI6").Value = -3.74
$ws.Range("AJ6").Value = 6138294

# Rows 7-9 (2019/12(E), 2020/12(E), 2021/12(E) forecast columns) had all data
# values removed, leaving only the A/B/C (index/period/name) columns intact.
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
